$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "27.207.70",
# "307.65") but must stay as plain text, exactly as they were authored.
# Force text formatting on the whole column first so assigning a
# numeric-looking string does not get reinterpreted as a number/date.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "27.201.81"

# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.904.95"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.31%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "307.66"
$ws.Range("E5").Value = "  +0.61%  "

# Row 6 (USDC)
$ws.Range("E6").Value = "  +0.28%  "

# Row 7 (XRP)
$ws.Range("D7").Value = "0.5269"
$ws.Range("E7").Value = "  +0.60%  "

# Row 9 (Dogecoin)
$ws.Range("D9").Value = "0.07301"
$ws.Range("E9").Value = "  +0.73%  "

# Row 10 (Solana)
$ws.Range("D10").Value = "21.57"
$ws.Range("E10").Value = "  +2.01%  "

# Row 11 (Polygon)
$ws.Range("D11").Value = "0.9050"
$ws.Range("E11").Value = "  +0.56%  "

# Row 12 (TRON)
$ws.Range("D12").Value = "0.08085"
$ws.Range("E12").Value = "  -4.15%  "

# Row 13 (Litecoin)
$ws.Range("D13").Value = "95.94"
$ws.Range("E13").Value = "  +0.90%  "

# Row 14 (Polkadot)
$ws.Range("D14").Value = "5.370"
$ws.Range("E14").Value = "  +1.51%  "

# Row 15 (WrappedEther)
$ws.Range("D15").Value = "1.795.62"
$ws.Range("E15").Value = "  -5.60%  "

# Row 16 (BinanceUSD)
$ws.Range("E16").Value = "  +0.31%  "

# Row 17 (ShibaInu)
$ws.Range("D17").Value = "0.000008676"
$ws.Range("E17").Value = "  +0.62%  "

# Row 18 (Avalanche)
$ws.Range("E18").Value = "  +1.17%  "

# Row 19 (Dai)
$ws.Range("E19").Value = "  +0.24%  "

# Row 20 (WrappedBTC)
$ws.Range("D20").Value = "27.244.01"

# Row 21 (Uniswap)
$ws.Range("D21").Value = "5.126"
$ws.Range("E21").Value = "  +1.18%  "

# Row 22 (Cosmos)
$ws.Range("E22").Value = "  +2.02%  "

# Row 23 (Chainlink)
$ws.Range("D23").Value = "6.493"
$ws.Range("E23").Value = "  +1.01%  "

# Row 24 (LidoDAOToken)
$ws.Range("D24").Value = "2.342"
$ws.Range("E24").Value = "  +2.62%  "

# Row 25 (Monero)
$ws.Range("D25").Value = "150.20"
$ws.Range("E25").Value = "  +2.19%  "

# Row 26 (EthereumClassic)
$ws.Range("D26").Value = "18.26"
$ws.Range("E26").Value = "  +0.44%  "

# Row 27 (Toncoin)
$ws.Range("D27").Value = "1.744"
$ws.Range("E27").Value = "  -0.37%  "

# Row 28 (BitcoinCash)
$ws.Range("D28").Value = "117.12"
$ws.Range("E28").Value = "  +1.87%  "

# Row 29 (InternetComputer(DFINITY))
$ws.Range("D29").Value = "4.844"
$ws.Range("E29").Value = "  +0.57%  "

# Row 30 (Filecoin)
$ws.Range("D30").Value = "4.871"
$ws.Range("E30").Value = "  -0.40%  "

# Row 31 (Stellar)
$ws.Range("D31").Value = "0.09232"
$ws.Range("E31").Value = "  -0.34%  "

# Row 32 (ImmutableX)
$ws.Range("D32").Value = "0.8096"
$ws.Range("E32").Value = "  +0.21%  "

# Row 33 (Hedera)
$ws.Range("E33").Value = "  +0.07%  "

# Row 34 (ARBITRUM)
$ws.Range("E34").Value = "  -0.81%  "

# Row 35 (HuobiToken)
$ws.Range("D35").Value = "2.985"
$ws.Range("E35").Value = "  +1.29%  "

# Row 36 (MXToken)
$ws.Range("D36").Value = "3.363"
$ws.Range("E36").Value = "  -2.01%  "

# Row 37 (RenderToken)
$ws.Range("D37").Value = "2.709"
$ws.Range("E37").Value = "  +3.13%  "

# Row 38 (TheSandbox)
$ws.Range("D38").Value = "0.5739"
$ws.Range("E38").Value = "  -0.05%  "

# Row 39 (VeChain)
$ws.Range("E39").Value = "  +0.11%  "

# Row 41 (Aptos)
$ws.Range("D41").Value = "8.980"
$ws.Range("E41").Value = "  -0.19%  "

# Row 42 (FraxShare)
$ws.Range("D42").Value = "6.624"
$ws.Range("E42").Value = "  -0.21%  "

# Row 43 (Quant)
$ws.Range("D43").Value = "116.71"
$ws.Range("E43").Value = "  +0.42%  "

# Row 44 (Algorand)
$ws.Range("D44").Value = "0.1521"
$ws.Range("E44").Value = "  +0.54%  "

# Row 45 (Decentraland)
$ws.Range("D45").Value = "0.4923"
$ws.Range("E45").Value = "  +1.09%  "

# Row 46 and 47 swap places: EnergySwap <-> PaxDollar
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.16"
$ws.Range("E47").Value = "  +0.19%  "

# Row 48 (NEARProtocol)
$ws.Range("D48").Value = "1.638"
$ws.Range("E48").Value = "  +1.63%  "

# Row 49 (Elrond)
$ws.Range("E49").Value = "  +2.94%  "

# Row 50 (Aave)
$ws.Range("D50").Value = "64.29"

# Row 51 (Cronos)
$ws.Range("D51").Value = "0.05963"
$ws.Range("E51").Value = "  +0.27%  "
